# Updated symbol list on Sun Dec 25 07:29:20 UTC 2022 with GitHub Actions
#
# This script applies the price/coin refresh to the "cryptos" worksheet.
# Column D ("Price") holds numeric-looking values that are stored as TEXT
# in the workbook (so leading/trailing zero formatting like "5.410" or
# "0.06009" is preserved exactly as scraped). We force those cells to a
# text number-format before writing so Excel does not silently convert
# them into real floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Price-only updates (Column D) ---------------------------------------
$priceUpdates = [ordered]@{
    "D2"  = "245.18"
    "D3"  = "23.05"
    "D4"  = "5.410"
    "D5"  = "0.06009"
    "D6"  = "3.389"
    "D7"  = "0.8087"
    "D8"  = "0.9332"
    "D9"  = "0.1421"
    "D10" = "0.07449"
    "D11" = "0.03368"
    "D12" = "0.03023"
    "D13" = "0.09355"
    "D14" = "3.949"
    "D15" = "0.001602"
    "D16" = "0.04827"
    "D40" = "0.03973"
    "D41" = "0.006343"
    "D42" = "0.1073"
    "D43" = "0.002901"
    "D44" = "0.006199"
    "D45" = "0.00005203"
    "D47" = "0.0005803"
    "D48" = "0.9804"
    "D49" = "0.002027"
}

foreach ($addr in $priceUpdates.Keys) {
    Set-TextValue $addr $priceUpdates[$addr]
}

# --- Coin list re-ranking (rows 17-24) ------------------------------------
# Coin, Link, Price, Volume(1h) shift up one rank; "One" moves to the
# bottom of this block (row 24) with a refreshed price/volume label.
$rowUpdates = @{
    17 = @("TigerCash",    "https://coinranking.com/coin/6hIn06L2+tigercash-tch",            "0.005469",   "16TigerCashTCH")
    18 = @("HotbitToken",  "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb",       "0.004151",   "17HotbitTokenHTB")
    19 = @("BitKan",       "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan",           "0.0009836",  "18BitKanKAN")
    20 = @("NitroEx",      "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx",            "0.00007103", "19NitroExNTX")
    21 = @("LEO",          "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo",               "3.662",      "20LEOLEO")
    22 = @("KuCoinToken",  "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs",      "6.433",      "21KuCoinTokenKCS")
    23 = @("BTSEToken",    "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse",        "2.187",      "22BTSETokenBTSE")
    24 = @("One",          "https://coinranking.com/coin/6Lga5NiXX3rT+one-one",               "0.01121",    "23OneONE")
}

foreach ($row in $rowUpdates.Keys) {
    $vals = $rowUpdates[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    Set-TextValue "D$row" $vals[2]
    $ws.Range("E$row").Value = $vals[3]
}
